# Regenerate orders with updated distance/size codes.
# The experiment's distance and size condition codes were renumbered:
#   D64 -> D69, D51 -> D55, D80 -> D86 (Distance column / embedded in labels & filenames)
#   S30 -> S31 (Size column / embedded in labels & filenames)
# These tokens appear throughout the "Condition", "Filename_Left",
# "Filename_Right", "Distance" and "Size" columns (e.g. "Face12_D64_S30",
# "Face12_D64_S30_l.png", "D64", "S30"), so apply the substitutions across
# every used cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("S30", "S31")
